# Add a "Total" row to the BOM sheet that sums the Cost column, format
# the Unit Cost / Cost columns as currency, and bold the Total row label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Currency number format used by the Unit Cost (B) and Cost (D) columns.
$currencyFormat = '"$"#,##0.00'

# Apply the currency number format to the existing Unit Cost / Cost data
# (columns B and D, rows 1-6) plus set the column-level default format so
# any new cell typed into those columns inherits it too.
$ws.Columns.Item("B").NumberFormat = $currencyFormat
$ws.Columns.Item("D").NumberFormat = $currencyFormat
$ws.Columns.Item("B").ColumnWidth = 8.7265625
$ws.Columns.Item("D").ColumnWidth = 8.7265625

# New "Total" row.
$ws.Range("A8").Value = "Total"
$ws.Range("A8").Font.Bold = $true

$ws.Range("D8").Formula = "=SUM(D2:D6)"
$ws.Range("D8").NumberFormat = $currencyFormat
$ws.Range("D8").Font.Bold = $true

# Update selection to match the author's final cursor position.
$ws.Range("C10").Select()

# Set up the page for printing (portrait orientation).
$ws.PageSetup.Orientation = 1

$wb.Save()
